$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the "PostalCode" column (currently column F),
# which pushes PostalCode to column G and leaves a blank column F.
$ws.Columns.Item(6).Insert()

# Give the new column its header text.
$ws.Cells.Item(1, 6).Value = "ManagingDirector"

# Match the new column's width to its neighbour (PhysicalAddress, column E).
$ws.Columns.Item(6).ColumnWidth = $ws.Columns.Item(5).ColumnWidth

# Copy the header formatting from the neighbouring header cell (E1) onto
# the new header cell (F1) so it matches the other headers.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Extend the table ("Table1") so it now covers the new column too.
$tbl = $ws.ListObjects.Item(1)
$tbl.ListColumns.Add() | Out-Null

# The table-extend above leaves the new (7th) header cell blank; restore
# its text (it is the PostalCode column that got shifted to G1).
$ws.Cells.Item(1, 7).Value = "PostalCode"

# Leave the selection where the author's edit left it.
$ws.Range("G5").Select()
